$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# COVID-19 data refresh: update per-country stats (columns B-H =
# Casos totales, Nuevos casos, Casos activos, Recuperados, Casos
# criticos, Muertes hoy, Muertes) for every row whose figures changed
# between the 16 Jul 23:31 snapshot and the 17 Jul 00:48 snapshot.

# Row 4: Estados Unidos
$ws.Range("B4").Value = 3680868
$ws.Range("C4").Value = 64041
$ws.Range("D4").Value = 1672403
$ws.Range("E4").Value = 1867511
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 810
$ws.Range("H4").Value = 140954

# Row 8: Peru
$ws.Range("B8").Value = 341586
$ws.Range("C8").Value = 3862
$ws.Range("D8").Value = 230994
$ws.Range("E8").Value = 97977
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 198
$ws.Range("H8").Value = 12615

# Row 22: Colombia
$ws.Range("B22").Value = 173206
$ws.Range("C22").Value = 8037
$ws.Range("D22").Value = 76164
$ws.Range("E22").Value = 91013
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 215
$ws.Range("H22").Value = 6029

# Row 23: Argentina
$ws.Range("B23").Value = 114783
$ws.Range("C23").Value = 3637
$ws.Range("D23").Value = 49120
$ws.Range("E23").Value = 63551
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 62
$ws.Range("H23").Value = 2112

# Row 59: Japon
$ws.Range("B59").Value = 22890
$ws.Range("C59").Value = 382
$ws.Range("D59").Value = 18814
$ws.Range("E59").Value = 3091
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 985

# Row 66: Camerun
$ws.Range("B66").Value = 16157
$ws.Range("C66").Value = 984
$ws.Range("D66").Value = 13728
$ws.Range("E66").Value = 2056
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 14
$ws.Range("H66").Value = 373

# Row 67: Uzbekistan
$ws.Range("B67").Value = 15066
$ws.Range("C67").Value = 485
$ws.Range("D67").Value = 8783
$ws.Range("E67").Value = 6208
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 4
$ws.Range("H67").Value = 75

# Row 79: Noruega
$ws.Range("B79").Value = 9015
$ws.Range("C79").Value = 4
$ws.Range("D79").Value = 8138
$ws.Range("E79").Value = 623
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 1
$ws.Range("H79").Value = 254

# Row 83: Consejo Danes para los Refugiados
$ws.Range("B83").Value = 8475
$ws.Range("C83").Value = 294
$ws.Range("D83").Value = 2430
$ws.Range("E83").Value = 5897
$ws.Range("F83").Value = 0
$ws.Range("G83").Value = 2
$ws.Range("H83").Value = 148

# Row 84: Etiopia
$ws.Range("B84").Value = 8199
$ws.Range("C84").Value = 36
$ws.Range("D84").Value = 4248
$ws.Range("E84").Value = 3758
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 193

# Row 85: Bulgaria
$ws.Range("B85").Value = 8144
$ws.Range("C85").Value = 267
$ws.Range("D85").Value = 3927
$ws.Range("E85").Value = 3924
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 4
$ws.Range("H85").Value = 293

# Row 90: Tayikistan
$ws.Range("B90").Value = 6741
$ws.Range("C90").Value = 46
$ws.Range("D90").Value = 5431
$ws.Range("E90").Value = 1254
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 56

# Row 97: Republica de Yibuti
$ws.Range("B97").Value = 4993
$ws.Range("C97").Value = 8
$ws.Range("D97").Value = 4796
$ws.Range("E97").Value = 141
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 56

# Row 103: Tailandia
$ws.Range("B103").Value = 3342
$ws.Range("C103").Value = 144
$ws.Range("D103").Value = 1379
$ws.Range("E103").Value = 1936
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 27

# Row 104: Paraguay
$ws.Range("B104").Value = 3236
$ws.Range("C104").Value = 4
$ws.Range("D104").Value = 3095
$ws.Range("E104").Value = 83
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 58

# Row 124: Cabo Verde
$ws.Range("B124").Value = 1894
$ws.Range("C124").Value = 114
$ws.Range("D124").Value = 902
$ws.Range("E124").Value = 973
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 19

# Row 134: Tunez
$ws.Range("B134").Value = 1362
$ws.Range("C134").Value = 273
$ws.Range("D134").Value = 425
$ws.Range("E134").Value = 914
$ws.Range("F134").Value = 0
$ws.Range("G134").Value = 3
$ws.Range("H134").Value = 23

# Row 135: Montenegro
$ws.Range("B135").Value = 1327
$ws.Range("C135").Value = 8
$ws.Range("D135").Value = 1093
$ws.Range("E135").Value = 184
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 50

# Row 136: Jordania
$ws.Range("B136").Value = 1287
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 330
$ws.Range("E136").Value = 933
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 24

# Row 137: Letonia
$ws.Range("B137").Value = 1206
$ws.Range("C137").Value = 5
$ws.Range("D137").Value = 1019
$ws.Range("E137").Value = 177
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 10

# Row 138: Niger
$ws.Range("B138").Value = 1179
$ws.Range("C138").Value = 1
$ws.Range("D138").Value = 1022
$ws.Range("E138").Value = 126
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 31

# Row 139: Zimbabue
$ws.Range("B139").Value = 1102
$ws.Range("C139").Value = 2
$ws.Range("D139").Value = 993
$ws.Range("E139").Value = 40
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 69

# Row 151: Togo
$ws.Range("B151").Value = 749
$ws.Range("C151").Value = 9
$ws.Range("D151").Value = 543
$ws.Range("E151").Value = 191
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 15

# Row 181: Bahamas
$ws.Range("B181").Value = 124
$ws.Range("C181").Value = 5
$ws.Range("D181").Value = 91
$ws.Range("E181").Value = 22
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 11

# Update the "last refreshed" timestamp string (row 1).
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 00:48"
